# Daily attendance processing - 2025-10-12 17:17:34
# Normalize the "Recorded By" column (G) so that multi-author entries are
# listed in reverse order (most-recent editor first).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val.ToString().Contains(",")) {
        $parts = $val.ToString().Split(",")
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        # Reverse the list order (manual loop - [array]::Reverse is not
        # reliable against this COM-interop shim).
        $reversed = @()
        for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
            $reversed += $trimmed[$i]
        }

        $cell.Value2 = [string]::Join(", ", $reversed)
    }
}
